$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 24.290427
$ws.Range("H2").Value = 72.871281
$ws.Range("I2").Value = 0.04034614150728122
$ws.Range("J2").Value = 0.04034614150728123
$ws.Range("M2").Value = 44.40220133333333
$ws.Range("N2").Value = 133.206604
$ws.Range("O2").Value = 0.9893265572082102
$ws.Range("P2").Value = 0.9893265572082101
$ws.Range("Q2").Value = 1078.548430126636
$ws.Range("R2").Value = 9706.935871139724
$ws.Range("S2").Value = 0.0399155092740338
$ws.Range("T2").Value = 0.0399155092740338

$ws.Range("G3").Value = 24.290427
$ws.Range("H3").Value = 72.871281
$ws.Range("I3").Value = 0.04034614150728122
$ws.Range("J3").Value = 0.04034614150728123
$ws.Range("M3").Value = 0.401961
$ws.Range("N3").Value = 1.205883
$ws.Range("O3").Value = 0.008956103083191794
$ws.Range("P3").Value = 0.008956103083191792
$ws.Range("Q3").Value = 9.763804327347
$ws.Range("R3").Value = 87.874238946123
$ws.Range("S3").Value = 0.0003613442023482538
$ws.Range("T3").Value = 0.0003613442023482537

$ws.Range("G4").Value = 24.290427
$ws.Range("H4").Value = 72.871281
$ws.Range("I4").Value = 0.04034614150728122
$ws.Range("J4").Value = 0.04034614150728123
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07707633333333333
$ws.Range("N4").Value = 0.231229
$ws.Range("O4").Value = 0.00171733970859806
$ws.Range("P4").Value = 0.00171733970859806
$ws.Range("Q4").Value = 1.872217048261
$ws.Range("R4").Value = 16.849953434349
$ws.Range("S4").Value = 0.00006928803089917043
$ws.Range("T4").Value = 0.00006928803089917045

$ws.Range("I5").Value = 0.9165963909899637
$ws.Range("J5").Value = 0.9165963909899638
$ws.Range("M5").Value = 44.40220133333333
$ws.Range("N5").Value = 133.206604
$ws.Range("O5").Value = 0.9893265572082102
$ws.Range("P5").Value = 0.9893265572082101
$ws.Range("Q5").Value = 24502.80402609392
$ws.Range("R5").Value = 220525.2362348453
$ws.Range("S5").Value = 0.9068131518475713
$ws.Range("T5").Value = 0.9068131518475713

$ws.Range("I6").Value = 0.9165963909899637
$ws.Range("J6").Value = 0.9165963909899638
$ws.Range("M6").Value = 0.401961
$ws.Range("N6").Value = 1.205883
$ws.Range("O6").Value = 0.008956103083191794
$ws.Range("P6").Value = 0.008956103083191792
$ws.Range("Q6").Value = 221.817191791769
$ws.Range("R6").Value = 1996.354726125921
$ws.Range("S6").Value = 0.008209131763387684
$ws.Range("T6").Value = 0.008209131763387684

$ws.Range("I7").Value = 0.9165963909899637
$ws.Range("J7").Value = 0.9165963909899638
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.07707633333333333
$ws.Range("N7").Value = 0.231229
$ws.Range("O7").Value = 0.00171733970859806
$ws.Range("P7").Value = 0.00171733970859806
$ws.Range("Q7").Value = 42.53361846946922
$ws.Range("R7").Value = 382.802566225223
$ws.Range("S7").Value = 0.001574107379004738
$ws.Range("T7").Value = 0.001574107379004738

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.741047
$ws.Range("H8").Value = 2.223141
$ws.Range("I8").Value = 0.001230871203933394
$ws.Range("J8").Value = 0.001230871203933394
$ws.Range("M8").Value = 44.40220133333333
$ws.Range("N8").Value = 133.206604
$ws.Range("O8").Value = 0.9893265572082102
$ws.Range("P8").Value = 0.9893265572082101
$ws.Range("Q8").Value = 32.90411809146266
$ws.Range("R8").Value = 296.137062823164
$ws.Range("S8").Value = 0.00121773357055415
$ws.Range("T8").Value = 0.00121773357055415

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.741047
$ws.Range("H9").Value = 2.223141
$ws.Range("I9").Value = 0.001230871203933394
$ws.Range("J9").Value = 0.001230871203933394
$ws.Range("M9").Value = 0.401961
$ws.Range("N9").Value = 1.205883
$ws.Range("O9").Value = 0.008956103083191794
$ws.Range("P9").Value = 0.008956103083191792
$ws.Range("Q9").Value = 0.297871993167
$ws.Range("R9").Value = 2.680847938503
$ws.Range("S9").Value = 0.00001102380938455987
$ws.Range("T9").Value = 0.00001102380938455987

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.741047
$ws.Range("H10").Value = 2.223141
$ws.Range("I10").Value = 0.001230871203933394
$ws.Range("J10").Value = 0.001230871203933394
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.07707633333333333
$ws.Range("N10").Value = 0.231229
$ws.Range("O10").Value = 0.00171733970859806
$ws.Range("P10").Value = 0.00171733970859806
$ws.Range("Q10").Value = 0.05711718558766667
$ws.Range("R10").Value = 0.514054670289
$ws.Range("S10").Value = 0.000002113823994684719
$ws.Range("T10").Value = 0.000002113823994684719

$ws.Range("G11").Value = 24.947199
$ws.Range("H11").Value = 74.84159699999999
$ws.Range("I11").Value = 0.04143703283043582
$ws.Range("J11").Value = 0.04143703283043582
$ws.Range("M11").Value = 44.40220133333333
$ws.Range("N11").Value = 133.206604
$ws.Range("O11").Value = 0.9893265572082102
$ws.Range("P11").Value = 0.9893265572082101
$ws.Range("Q11").Value = 1107.710552700732
$ws.Range("R11").Value = 9969.394974306588
$ws.Range("S11").Value = 0.04099475703105864
$ws.Range("T11").Value = 0.04099475703105863

$ws.Range("G12").Value = 24.947199
$ws.Range("H12").Value = 74.84159699999999
$ws.Range("I12").Value = 0.04143703283043582
$ws.Range("J12").Value = 0.04143703283043582
$ws.Range("M12").Value = 0.401961
$ws.Range("N12").Value = 1.205883
$ws.Range("O12").Value = 0.008956103083191794
$ws.Range("P12").Value = 0.008956103083191792
$ws.Range("Q12").Value = 10.027801057239
$ws.Range("R12").Value = 90.250209515151
$ws.Range("S12").Value = 0.0003711143374909858
$ws.Range("T12").Value = 0.0003711143374909857

$ws.Range("G13").Value = 24.947199
$ws.Range("H13").Value = 74.84159699999999
$ws.Range("I13").Value = 0.04143703283043582
$ws.Range("J13").Value = 0.04143703283043582
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.07707633333333333
$ws.Range("N13").Value = 0.231229
$ws.Range("O13").Value = 0.00171733970859806
$ws.Range("P13").Value = 0.00171733970859806
$ws.Range("Q13").Value = 1.922838625857
$ws.Range("R13").Value = 17.305547632713
$ws.Range("S13").Value = 0.00007116146188618892
$ws.Range("T13").Value = 0.0000711614618861889

$ws.Range("G14").Value = 0.234537
$ws.Range("H14").Value = 0.7036110000000001
$ws.Range("I14").Value = 0.0003895634683858467
$ws.Range("J14").Value = 0.0003895634683858468
$ws.Range("M14").Value = 44.40220133333333
$ws.Range("N14").Value = 133.206604
$ws.Range("O14").Value = 0.9893265572082102
$ws.Range("P14").Value = 0.9893265572082101
$ws.Range("Q14").Value = 10.413959094116
$ws.Range("R14").Value = 93.72563184704401
$ws.Range("S14").Value = 0.0003854054849922591
$ws.Range("T14").Value = 0.0003854054849922591

$ws.Range("G15").Value = 0.234537
$ws.Range("H15").Value = 0.7036110000000001
$ws.Range("I15").Value = 0.0003895634683858467
$ws.Range("J15").Value = 0.0003895634683858468
$ws.Range("M15").Value = 0.401961
$ws.Range("N15").Value = 1.205883
$ws.Range("O15").Value = 0.008956103083191794
$ws.Range("P15").Value = 0.008956103083191792
$ws.Range("Q15").Value = 0.09427472705700002
$ws.Range("R15").Value = 0.8484725435130002
$ws.Range("S15").Value = 0.000003488970580309371
$ws.Range("T15").Value = 0.00000348897058030937

$ws.Range("G16").Value = 0.234537
$ws.Range("H16").Value = 0.7036110000000001
$ws.Range("I16").Value = 0.0003895634683858467
$ws.Range("J16").Value = 0.0003895634683858468
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.07707633333333333
$ws.Range("N16").Value = 0.231229
$ws.Range("O16").Value = 0.00171733970859806
$ws.Range("P16").Value = 0.00171733970859806
$ws.Range("Q16").Value = 0.018077251991
$ws.Range("R16").Value = 0.162695267919
$ws.Range("S16").Value = 0.0000006690128132781997
$ws.Range("T16").Value = 0.0000006690128132781997
